# Applies the 2026-02-18 05:50 automatic-update diff to Dades_Meteo.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-18 05:48:35'
$ws.Range('H2').NumberFormat = "@"
$ws.Range('H2').Value = '71%'
$ws.Range('N2').Value = '-1.9 °C 5:19 TU'
$ws.Range('O2').Value = '-1.0 °C'
$ws.Range('E3').Value = '2026-02-18 05:48:38'
$ws.Range('E4').Value = '2026-02-18 05:48:40'
$ws.Range('H4').NumberFormat = "@"
$ws.Range('H4').Value = '97%'
$ws.Range('J4').Value = '1018.5 hPa'
$ws.Range('N4').Value = '5.1 °C 5:12 TU'
$ws.Range('O4').Value = '7.3 °C'
$ws.Range('E5').Value = '2026-02-18 05:48:42'
$ws.Range('H5').NumberFormat = "@"
$ws.Range('H5').Value = '89%'
$ws.Range('M5').Value = '-0.4 °C 5:10 TU'
$ws.Range('O5').Value = '-2.1 °C'
$ws.Range('E6').Value = '2026-02-18 05:48:45'
$ws.Range('H6').NumberFormat = "@"
$ws.Range('H6').Value = '96%'
$ws.Range('J6').Value = '1018.2 hPa'
$ws.Range('N6').Value = '6.0 °C 5:20 TU'
$ws.Range('O6').Value = '8.1 °C'
$ws.Range('E7').Value = '2026-02-18 05:48:47'
$ws.Range('J7').Value = '1018.4 hPa'
$ws.Range('N7').Value = '11.7 °C 5:29 TU'
$ws.Range('E8').Value = '2026-02-18 05:48:50'
$ws.Range('J8').Value = '1018.5 hPa'
$ws.Range('N8').Value = '8.0 °C 5:03 TU'
$ws.Range('O8').Value = '8.8 °C'
$ws.Range('E9').Value = '2026-02-18 05:48:52'
$ws.Range('N9').Value = '4.1 °C 5:22 TU'
$ws.Range('O9').Value = '5.0 °C'
$ws.Range('E10').Value = '2026-02-18 05:48:54'
$ws.Range('N10').Value = '4.2 °C 5:22 TU'
$ws.Range('O10').Value = '7.0 °C'
$ws.Range('E11').Value = '2026-02-18 05:48:57'
$ws.Range('H11').NumberFormat = "@"
$ws.Range('H11').Value = '98%'
$ws.Range('N11').Value = '-0.1 °C 5:23 TU'
$ws.Range('O11').Value = '2.1 °C'
$ws.Range('E12').Value = '2026-02-18 05:48:59'
$ws.Range('N12').Value = '4.5 °C 5:29 TU'
$ws.Range('O12').Value = '6.0 °C'
$ws.Range('E13').Value = '2026-02-18 05:49:01'
$ws.Range('O13').Value = '-2.2 °C'
$ws.Range('E14').Value = '2026-02-18 05:49:03'
$ws.Range('L14').Value = '11.2 km/h - 304º 5:07 TU'
$ws.Range('O14').Value = '10.1 °C'
$ws.Range('E15').Value = '2026-02-18 05:49:06'
$ws.Range('E16').Value = '2026-02-18 05:49:08'
$ws.Range('H16').NumberFormat = "@"
$ws.Range('H16').Value = '37%'
$ws.Range('M16').Value = '1.5 °C 5:17 TU'
$ws.Range('O16').Value = '0.7 °C'
$ws.Range('E17').Value = '2026-02-18 05:49:10'
$ws.Range('H17').NumberFormat = "@"
$ws.Range('H17').Value = '89%'
$ws.Range('L17').Value = '49.7 km/h - 283º 5:26 TU'
$ws.Range('M17').Value = '3.8 °C 5:28 TU'
$ws.Range('O17').Value = '2.1 °C'
$ws.Range('E18').Value = '2026-02-18 05:49:12'
$ws.Range('H18').NumberFormat = "@"
$ws.Range('H18').Value = '95%'
$ws.Range('J18').Value = '1018.5 hPa'
$ws.Range('N18').Value = '5.0 °C 5:16 TU'
$ws.Range('O18').Value = '7.6 °C'
$ws.Range('E19').Value = '2026-02-18 05:49:15'
$ws.Range('N19').Value = '5.3 °C 5:24 TU'
$ws.Range('E20').Value = '2026-02-18 05:49:17'
$ws.Range('H20').NumberFormat = "@"
$ws.Range('H20').Value = '80%'
$ws.Range('K20').Value = '-0.1 MJ/m2'
$ws.Range('E21').Value = '2026-02-18 05:49:20'
$ws.Range('H21').NumberFormat = "@"
$ws.Range('H21').Value = '84%'
$ws.Range('O21').Value = '2.0 °C'
$ws.Range('E22').Value = '2026-02-18 05:49:22'
$ws.Range('E23').Value = '2026-02-18 05:49:24'
$ws.Range('O23').Value = '0.4 °C'
$ws.Range('E24').Value = '2026-02-18 05:49:27'
$ws.Range('N24').Value = '2.8 °C 5:29 TU'
$ws.Range('O24').Value = '5.3 °C'
$ws.Range('E25').Value = '2026-02-18 05:49:29'
$ws.Range('M25').Value = '1.2 °C 5:25 TU'
$ws.Range('O25').Value = '-0.4 °C'
$ws.Range('E26').Value = '2026-02-18 05:49:31'
$ws.Range('E27').Value = '2026-02-18 05:49:33'
$ws.Range('G27').Value = '168 cm'
$ws.Range('H27').NumberFormat = "@"
$ws.Range('H27').Value = '53%'
$ws.Range('M27').Value = '2.0 °C 5:28 TU'
$ws.Range('O27').Value = '0.8 °C'
$ws.Range('E28').Value = '2026-02-18 05:49:36'
$ws.Range('J28').Value = '1018.9 hPa'
$ws.Range('O28').Value = '4.9 °C'
$ws.Range('E29').Value = '2026-02-18 05:49:38'
$ws.Range('H29').NumberFormat = "@"
$ws.Range('H29').Value = '91%'
$ws.Range('N29').Value = '7.2 °C 5:04 TU'
$ws.Range('O29').Value = '9.5 °C'
$ws.Range('E30').Value = '2026-02-18 05:49:41'
$ws.Range('J30').Value = '1018.5 hPa'
$ws.Range('N30').Value = '5.2 °C 5:29 TU'
$ws.Range('E31').Value = '2026-02-18 05:49:43'
$ws.Range('H31').NumberFormat = "@"
$ws.Range('H31').Value = '77%'
$ws.Range('J31').Value = '1017.1 hPa'
$ws.Range('N31').Value = '9.8 °C 5:18 TU'
$ws.Range('E32').Value = '2026-02-18 05:49:45'
$ws.Range('L32').Value = '12.6 km/h - 293º 5:27 TU'
$ws.Range('M32').Value = '5.9 °C 5:12 TU'
$ws.Range('E33').Value = '2026-02-18 05:49:48'
$ws.Range('J33').Value = '1021.3 hPa'
$ws.Range('N33').Value = '-1.5 °C 5:27 TU'
$ws.Range('O33').Value = '0.0 °C'
$ws.Range('E34').Value = '2026-02-18 05:49:50'
$ws.Range('H34').NumberFormat = "@"
$ws.Range('H34').Value = '54%'
$ws.Range('M34').Value = '3.6 °C 5:04 TU'
$ws.Range('O34').Value = '-0.1 °C'
$ws.Range('E35').Value = '2026-02-18 05:49:52'
$ws.Range('H35').NumberFormat = "@"
$ws.Range('H35').Value = '81%'
$ws.Range('J35').Value = '1019.1 hPa'
$ws.Range('N35').Value = '5.0 °C 5:21 TU'
$ws.Range('O35').Value = '7.3 °C'
$ws.Range('E36').Value = '2026-02-18 05:49:55'
$ws.Range('J36').Value = '1018.4 hPa'
$ws.Range('N36').Value = '5.6 °C 5:00 TU'
$ws.Range('O36').Value = '8.3 °C'
$ws.Range('E37').Value = '2026-02-18 05:49:57'
$ws.Range('J37').Value = '1021.3 hPa'
$ws.Range('L37').Value = '9.7 km/h - 353º 5:16 TU'
$ws.Range('N37').Value = '0.2 °C 5:12 TU'
$ws.Range('O37').Value = '1.4 °C'
$ws.Range('E38').Value = '2026-02-18 05:50:00'
$ws.Range('H38').NumberFormat = "@"
$ws.Range('H38').Value = '93%'
$ws.Range('N38').Value = '5.9 °C 5:08 TU'
$ws.Range('O38').Value = '9.0 °C'
$ws.Range('E39').Value = '2026-02-18 05:50:02'
$ws.Range('E40').Value = '2026-02-18 05:50:04'
$ws.Range('O40').Value = '0.8 °C'
$ws.Range('E41').Value = '2026-02-18 05:50:07'
$ws.Range('J41').Value = '1018.1 hPa'
$ws.Range('N41').Value = '6.9 °C 5:12 TU'
$ws.Range('O41').Value = '8.3 °C'
$ws.Range('E42').Value = '2026-02-18 05:50:09'
$ws.Range('O42').Value = '8.3 °C'
$ws.Range('E43').Value = '2026-02-18 05:50:11'
$ws.Range('N43').Value = '5.9 °C 5:29 TU'
$ws.Range('O43').Value = '7.1 °C'
$ws.Range('E44').Value = '2026-02-18 05:50:14'
$ws.Range('H44').NumberFormat = "@"
$ws.Range('H44').Value = '67%'
$ws.Range('O44').Value = '-3.6 °C'
$ws.Range('E45').Value = '2026-02-18 05:50:16'
$ws.Range('J45').Value = '1020.7 hPa'
$ws.Range('E46').Value = '2026-02-18 05:50:19'
$ws.Range('J46').Value = '1019.0 hPa'
$ws.Range('N46').Value = '4.2 °C 5:17 TU'
$ws.Range('O46').Value = '6.3 °C'

Write-Host "Applied 146 cell updates."
